# Regenerate merged AHB files: rename diff-header columns and turn the
# sheet into a frozen-header Table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AHB-Diff")

# A1:J1 -> "..._old" becomes "..._FV2404"
$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

# L1:U1 -> "..._new" becomes "..._FV2410"
$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $fv2404Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2404Headers[$i]
}
# column K (11) stays "diff"
for ($i = 0; $i -lt $fv2410Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2410Headers[$i]
}

# Freeze the header row.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Turn the used range into a table.
$rng = $ws.Range("A1:U89")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = $null
